$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the report header (From/To timestamp) text in A1
$ws.Range("A1").Value = "From: 28/01/2018 at 2014`nTo: 28/01/2018 at 2153"

# Update the ping-success-percentage data cells
$ws.Range("A6").Value = 99.95
$ws.Range("B6").Value = 99.95
$ws.Range("C6").Value = 99.95
$ws.Range("D6").Value = 99.95
$ws.Range("E6").Value = 99.95
$ws.Range("F6").Value = 99.95
$ws.Range("A12").Value = 99.95
$ws.Range("B12").Value = 99.95
$ws.Range("C12").Value = 99.95
$ws.Range("D12").Value = 99.95
$ws.Range("E12").Value = 99.95
$ws.Range("F12").Value = 99.95
$ws.Range("A18").Value = 99.95
$ws.Range("B18").Value = 99.95
$ws.Range("C18").Value = 99.95
$ws.Range("D18").Value = 99.95
$ws.Range("E18").Value = 99.95
$ws.Range("A24").Value = 99.95
$ws.Range("B24").Value = 99.95
$ws.Range("C24").Value = 99.95
$ws.Range("D24").Value = 99.95
$ws.Range("E24").Value = 99.95
$ws.Range("F24").Value = 99.85
$ws.Range("A34").Value = 99.95
$ws.Range("B34").Value = 99.95
$ws.Range("C34").Value = 99.95
$ws.Range("D34").Value = 99.95
$ws.Range("E34").Value = 99.95
$ws.Range("F34").Value = 99.85
$ws.Range("A38").Value = 99.85
$ws.Range("B38").Value = 99.95
$ws.Range("C38").Value = 99.95
$ws.Range("D38").Value = 99.95
$ws.Range("A44").Value = 99.95
$ws.Range("B44").Value = 99.95
$ws.Range("C44").Value = 99.95
$ws.Range("D44").Value = 99.95
$ws.Range("E44").Value = 99.95
$ws.Range("F44").Value = 99.85
$ws.Range("A48").Value = 99.85
$ws.Range("B48").Value = 99.95
$ws.Range("C48").Value = 99.95
$ws.Range("D48").Value = 99.9
$ws.Range("A58").Value = 99.95
$ws.Range("B58").Value = 99.95
$ws.Range("C58").Value = 99.95
$ws.Range("D58").Value = 99.95
$ws.Range("E58").Value = 99.95
$ws.Range("F58").Value = 99.85
$ws.Range("A62").Value = 99.85
$ws.Range("B62").Value = 99.95
$ws.Range("C62").Value = 99.95
$ws.Range("D62").Value = 99.9
$ws.Range("A72").Value = 99.95
$ws.Range("B72").Value = 99.95
$ws.Range("C72").Value = 99.95
$ws.Range("D72").Value = 99.95
$ws.Range("E72").Value = 99.95
$ws.Range("F72").Value = 99.85
$ws.Range("A76").Value = 99.85
$ws.Range("B76").Value = 99.95
$ws.Range("C76").Value = 99.95
$ws.Range("D76").Value = 99.9
$ws.Range("E76").Value = 99.9
$ws.Range("F76").Value = 99.85
$ws.Range("A80").Value = 99.85
$ws.Range("B80").Value = 99.9
$ws.Range("C80").Value = 99.9
$ws.Range("D80").Value = 99.95
$ws.Range("A86").Value = 99.95
$ws.Range("B86").Value = 99.95
$ws.Range("C86").Value = 99.95
$ws.Range("D86").Value = 99.95
$ws.Range("E86").Value = 99.95
$ws.Range("F86").Value = 99.85
$ws.Range("A90").Value = 99.85
$ws.Range("B90").Value = 99.95
$ws.Range("C90").Value = 99.95
$ws.Range("D90").Value = 99.9
$ws.Range("E90").Value = 99.9
$ws.Range("F90").Value = 99.85
$ws.Range("A94").Value = 99.85
$ws.Range("B94").Value = 99.9
$ws.Range("C94").Value = 99.9
$ws.Range("D94").Value = 99.95
$ws.Range("E94").Value = 99.95
$ws.Range("F94").Value = 96.61
$ws.Range("A98").Value = 96.11
$ws.Range("B98").Value = 85.2
$ws.Range("A104").Value = 99.95
$ws.Range("B104").Value = 99.95
$ws.Range("C104").Value = 99.95
$ws.Range("D104").Value = 99.95
$ws.Range("E104").Value = 99.95
$ws.Range("F104").Value = 99.85
$ws.Range("A108").Value = 99.85
$ws.Range("B108").Value = 99.95
$ws.Range("C108").Value = 99.95
$ws.Range("D108").Value = 99.9
$ws.Range("E108").Value = 99.9
$ws.Range("F108").Value = 99.85
$ws.Range("A112").Value = 99.85
$ws.Range("B112").Value = 99.9
$ws.Range("C112").Value = 99.9
$ws.Range("D112").Value = 99.95
$ws.Range("E112").Value = 99.95
$ws.Range("F112").Value = 96.61
$ws.Range("A116").Value = 96.11
$ws.Range("B116").Value = 85.2
$ws.Range("C116").Value = 84.65
$ws.Range("D116").Value = 84.7
$ws.Range("A122").Value = 99.95
$ws.Range("B122").Value = 99.95
$ws.Range("C122").Value = 99.95
$ws.Range("D122").Value = 99.95
$ws.Range("E122").Value = 99.95
$ws.Range("F122").Value = 99.85
$ws.Range("A126").Value = 99.85
$ws.Range("B126").Value = 99.95
$ws.Range("C126").Value = 99.95
$ws.Range("D126").Value = 99.9
$ws.Range("E126").Value = 99.9
$ws.Range("F126").Value = 99.85
$ws.Range("A130").Value = 99.85
$ws.Range("B130").Value = 99.9
$ws.Range("C130").Value = 99.9
$ws.Range("D130").Value = 99.95
$ws.Range("E130").Value = 99.95
$ws.Range("F130").Value = 96.61
$ws.Range("A134").Value = 96.11
$ws.Range("B134").Value = 85.2
$ws.Range("C134").Value = 40.909999999999997
$ws.Range("D134").Value = 33.729999999999997
$ws.Range("E134").Value = 34.53
$ws.Range("F134").Value = 34.33
$ws.Range("A140").Value = 99.95
$ws.Range("B140").Value = 99.95
$ws.Range("C140").Value = 99.95
$ws.Range("D140").Value = 99.95
$ws.Range("E140").Value = 99.95
$ws.Range("F140").Value = 99.85
$ws.Range("A144").Value = 99.85
$ws.Range("B144").Value = 99.95
$ws.Range("C144").Value = 99.95
$ws.Range("D144").Value = 99.9
$ws.Range("E144").Value = 99.9
$ws.Range("F144").Value = 99.85
$ws.Range("A148").Value = 99.85
$ws.Range("B148").Value = 99.9
$ws.Range("C148").Value = 99.9
$ws.Range("D148").Value = 99.95
$ws.Range("E148").Value = 99.95
$ws.Range("F148").Value = 96.61
$ws.Range("A152").Value = 96.11
$ws.Range("B152").Value = 85.2
$ws.Range("C152").Value = 40.909999999999997
$ws.Range("D152").Value = 33.729999999999997

# Clear the explicit paper size from the page setup (paperSize attribute removed)
$ws.PageSetup.PaperSize = $null

